$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 6 - INODOROS
$ws.Range("D6").Value = -169.19
$ws.Range("E6").Value = 1076.356108615601
$ws.Range("F6").Value = -0.1865038810347487

# Row 10 - PANELES DECORATIVOS
$ws.Range("D10").Value = -213.19
$ws.Range("E10").Value = 601.297983534392
$ws.Range("F10").Value = -0.5493058866208771

# Row 12 - PORCELANATO
$ws.Range("D12").Value = 22586.98
$ws.Range("E12").Value = 5368
$ws.Range("F12").Value = 0.8079769686832186

# Row 14 - TOTAL
$ws.Range("D14").Value = 22424.59
$ws.Range("E14").Value = 19778.79110009468
$ws.Range("F14").Value = 0.5313458167442818
